$d = $word.ActiveDocument

# The document contains 7 occurrences of an "<id>...</id>" marker, each
# split across three runs: "<id>" (Courier New styled), the bare id text
# "p065v_N" (plain/simple styled), and "</id>" (Courier New styled).
# Collapse each triple into a single run carrying the full
# "<id>p065v_N</id>" text, using the first run's (Courier New) formatting,
# matching the target edit.

for ($i = 1; $i -le 7; $i++) {
    $old = "<id>p065v_$i</id>"
    $new = "<id>p065v_$i</id>"
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: pattern not found/replaced for index $i"
    }
}
